# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de sheets to reflect the latest report run.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-24 16:22:17"
$wsZh.Range("E3").Value = "2016-03-24 16:22:17"
$wsZh.Range("H2").Value = "2016-03-24 16:22:45"
$wsZh.Range("H3").Value = "2016-03-24 16:22:45"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-24 16:22:22"
$wsDe.Range("E3").Value = "2016-03-24 16:22:22"
$wsDe.Range("H2").Value = "2016-03-24 16:22:52"
$wsDe.Range("H3").Value = "2016-03-24 16:22:52"
